# CRE20-015-02 Special Support Scheme for Continuing Medical Care to HA
# Patients in GD (SSSCMC) - Phase 2 (Reimbursement)
#
# The "02" sheet (eHS(S)M0014-02: payment data file by transaction) had an
# extra "二级科室" (Secondary Department) header column that doesn't belong
# in the Phase 2 reimbursement template. Remove it by deleting the entire
# column D, which shifts every later header/cell one column to the left.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("02")
$ws.Columns("D").Delete()
